# Updated NATMI TPM-derived LR-pair stats for F2 -> F2rl1 (YoungD0).
# Sender-dependent metrics (E-J) depend only on sending cluster; target-dependent
# metrics (K-P) depend only on target cluster; edge metrics (Q-T) = sender * target.
# The update also adds the missing MuSCs-sending rows, completing the full 3x3
# (ECs/FAPs/MuSCs) sender x target cartesian product (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @(2, @("ECs", "F2", "F2rl1", "ECs", 3, 1, 0.3845463333333334, 1.153639, 0.1984850200147207, 0.1984850200147207, 1, 0.3333333333333333, 0.05881233333333333, 0.176437, 0.005710212285055511, 0.005710212285055511, 0.02261606713811111, 0.203544604243, 0.001133391599687547, 0.001133391599687547)),
    @(3, @("ECs", "F2", "F2rl1", "FAPs", 3, 1, 0.3845463333333334, 1.153639, 0.1984850200147207, 0.1984850200147207, 3, 1, 10.20548733333333, 30.616462, 0.990872081464405, 0.990872081464405, 3.924482733913111, 35.320344605218, 0.1966732649214904, 0.1966732649214904)),
    @(4, @("ECs", "F2", "F2rl1", "MuSCs", 3, 1, 0.3845463333333334, 1.153639, 0.1984850200147207, 0.1984850200147207, 1, 0.3333333333333333, 0.03520066666666666, 0.105602, 0.003417706250539468, 0.003417706250539468, 0.01353628729755556, 0.121826585678, 0.0006783634935427625, 0.0006783634935427624)),
    @(5, @("FAPs", "F2", "F2rl1", "ECs", 2, 0.6666666666666666, 1.110828, 3.332484, 0.5733580031870772, 0.5733580031870772, 1, 0.3333333333333333, 0.05881233333333333, 0.176437, 0.005710212285055511, 0.005710212285055511, 0.065330386612, 0.587973479508, 0.003273995913533745, 0.003273995913533745)),
    @(6, @("FAPs", "F2", "F2rl1", "FAPs", 2, 0.6666666666666666, 1.110828, 3.332484, 0.5733580031870772, 0.5733580031870772, 3, 1, 10.20548733333333, 30.616462, 0.990872081464405, 0.990872081464405, 11.336541083512, 102.028869751608, 0.5681244380422541, 0.5681244380422541)),
    @(7, @("FAPs", "F2", "F2rl1", "MuSCs", 2, 0.6666666666666666, 1.110828, 3.332484, 0.5733580031870772, 0.5733580031870772, 1, 0.3333333333333333, 0.03520066666666666, 0.105602, 0.003417706250539468, 0.003417706250539468, 0.039101886152, 0.351916975368, 0.001959569231289302, 0.001959569231289302)),
    @(8, @("MuSCs", "F2", "F2rl1", "ECs", 3, 1, 0.442033, 1.326099, 0.2281569767982021, 0.2281569767982021, 1, 0.3333333333333333, 0.05881233333333333, 0.176437, 0.005710212285055511, 0.005710212285055511, 0.02599699214033333, 0.233972929263, 0.001302824771834219, 0.001302824771834219)),
    @(9, @("MuSCs", "F2", "F2rl1", "FAPs", 3, 1, 0.442033, 1.326099, 0.2281569767982021, 0.2281569767982021, 3, 1, 10.20548733333333, 30.616462, 0.990872081464405, 0.990872081464405, 4.511162182415332, 40.600459641738, 0.2260743785006605, 0.2260743785006605)),
    @(10, @("MuSCs", "F2", "F2rl1", "MuSCs", 3, 1, 0.442033, 1.326099, 0.2281569767982021, 0.2281569767982021, 1, 0.3333333333333333, 0.03520066666666666, 0.105602, 0.003417706250539468, 0.003417706250539468, 0.01555985628866666, 0.140038706598, 0.0007797735257074038, 0.0007797735257074037))
)

foreach ($entry in $rowData) {
    $rowNum = $entry[0]
    $vals = $entry[1]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $vals[$c]
    }
}

Write-Output "edit applied"
